$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''28.223.55'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.29%  '
$ws.Range('D3').Value = '''1.919.89'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.99%  '
$ws.Range('E4').Value = '  -1.55%  '
$ws.Range('D5').Value = '''316.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.96%  '
$ws.Range('E6').Value = '  -1.39%  '
$ws.Range('D7').Value = '''0.4849'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.89%  '
$ws.Range('D8').Value = '''0.3851'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.15%  '
$ws.Range('D9').Value = '''0.07417'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('D10').Value = '''0.9489'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.14%  '
$ws.Range('D11').Value = '''20.94'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('D12').Value = '''0.07795'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.22%  '
$ws.Range('D13').Value = '''1.930.11'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.21%  '
$ws.Range('D14').Value = '''5.557'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.16%  '
$ws.Range('D15').Value = '''6.674'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').Value = '''92.11'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.91%  '
$ws.Range('E17').Value = '  -1.47%  '
$ws.Range('D18').Value = '''0.000008858'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('E19').Value = '  -1.37%  '
$ws.Range('D20').Value = '''28.226.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.12%  '
$ws.Range('D21').Value = '''15.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.40%  '
$ws.Range('D22').Value = '''5.174'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.69%  '
$ws.Range('D23').Value = '''2.169.44'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.97%  '
$ws.Range('D24').Value = '''10.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.28%  '
$ws.Range('D25').Value = '''1.932'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E26').Value = '  +1.34%  '
$ws.Range('D27').Value = '''18.69'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('D28').Value = '''2.108'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.85%  '
$ws.Range('D29').Value = '''117.35'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.08%  '
$ws.Range('D30').Value = '''5.016'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').Value = '''0.08908'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').Value = '''3.367'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('E33').Value = '  +4.63%  '
$ws.Range('D34').Value = '''0.7774'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.40%  '
$ws.Range('E35').Value = '  +2.53%  '
$ws.Range('D36').Value = '''2.769'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.99%  '
$ws.Range('D37').Value = '''0.02056'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('D38').Value = '''1.131'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '''0.05380'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.43%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '''0.5589'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.73%  '
$ws.Range('D41').Value = '''3.042'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.22%  '
$ws.Range('D42').Value = '''7.113'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').Value = '''8.579'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.10%  '
$ws.Range('D44').Value = '''0.1536'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').Value = '''10.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.55%  '
$ws.Range('E46').Value = '  +1.84%  '
$ws.Range('D47').Value = '''107.32'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.75%  '
$ws.Range('D48').Value = '''1.006'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.47%  '
$ws.Range('D49').Value = '''1.681'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').Value = '''69.49'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.98%  '
$ws.Range('D51').Value = '''0.06152'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.85%  '
